$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Pre_ISI (column I) value was normalized to 3.01
$rows = @(3,7,11,13,14,18,19,21,25,31,36,42,57,58,62,65,68,72,74,77,84,87,90,91,92,94)
foreach ($r in $rows) {
    $ws.Range("I$r").Value = 3.01
}

# Selection moved from L14 to P7
$ws.Range("P7").Select()
